# "Fixed StudyComb for Faceted Filters ICDC"
#
# The workbook has a single sheet ("startup") that stores a small table of
# Neo4j/Cypher queries used by a test-automation tool. Column D
# ("cartQuery") holds the same query text for the CasesTab/SamplesTab/
# FilesTab rows (D2:D4) - it is rewritten here with the corrected,
# much shorter "StudyComb" query.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$newQuery = "MATCH (demo:demographic)`n" +
            "WHERE demo.breed IN ['Yorkshire Terrier']`n" +
            "MATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`n" +
            "OPTIONAL MATCH (c)<-[*]-(samp:sample)`n" +
            "OPTIONAL MATCH (c)<-[*]-(f:file)`n" +
            "RETURN `n" +
            "`tcount(DISTINCT(f)) as number_of_files, `n" +
            "`tcount(DISTINCT(samp)) as number_of_sample, `n" +
            "`tcount(DISTINCT(c)) as number_of_cases, `n" +
            "`tcount(DISTINCT(s)) as number_of_study"

$ws.Range("D2:D4").Value = $newQuery

# The new query text wraps onto far fewer lines than the old one, so the
# (wrap-text, auto-fitted) rows shrink accordingly.
$ws.Rows("2:4").RowHeight = 244.8

# Reflect the author's final view state: zoomed in further and the
# selection left on C2 instead of the old C14.
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C2").Select() | Out-Null
